# Normalize the "Recorded By" (column G) cells so that whenever the list of
# recorders ends with "System" (e.g. "dnasr281@gmail.com, System" or
# "system, backup@backdoor.com, System"), "System" is moved to the front of
# the comma separated list, e.g. "System, dnasr281@gmail.com" /
# "System, system, backup@backdoor.com". Entries that do not end in
# ", System" (already start with System, are "System" alone, or don't
# mention System at all) are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$rowCount = $used.Rows.Count

for ($i = 1; $i -le $rowCount; $i++) {
    $cell = $ws.Cells.Item($i, 7)
    $orig = $cell.Text

    if ($orig -ne $null -and $orig.EndsWith(", System")) {
        $parts = $orig.Split(",")
        $trimmed = @()
        foreach ($p in $parts) {
            $trimmed += $p.Trim()
        }

        $n = $trimmed.Length
        $rest = $trimmed[0..($n - 2)]
        $newValue = "System, " + [string]::Join(", ", $rest)

        $cell.Value = $newValue
    }
}
